$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(42613.757210648146, 40),
    @(42613.884375000001, 42),
    @(42614.883136574077, 3),
    @(42615.884050925924, 69)
)

$row = 12
foreach ($entry in $data) {
    $date = $entry[0]
    $b = $entry[1]

    $ws.Cells.Item($row, 1).Value = $date
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = 0
    $ws.Cells.Item($row, 4).Value = 0
    $ws.Cells.Item($row, 5).Value = 0
    $ws.Cells.Item($row, 6).Value = 0
    $ws.Cells.Item($row, 7).Value = 0
    $ws.Cells.Item($row, 8).Value = 0
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 0
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = "Random"

    $row = $row + 1
}
